# Update the Sprint 2 Burndown Chart workbook:
# The sprint is now complete, so the "Actual" work-remaining value for the
# last day (day 7, cell C9) is recorded as 0 (previously blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("C9").Value = 0
